$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.63289999999998

$ws.Range("A7").Value = -22.02970000000001
$ws.Range("B7").Value = 4.515900000000001

$ws.Range("B15").Value = 4.992999999999997

$ws.Range("A16").Value = -21.97890000000002
$ws.Range("E16").Value = 16.443

$ws.Range("E19").Value = 16.3484

$ws.Range("B21").Value = 10.1033

$ws.Range("B22").Value = 9.916199999999998

$ws.Range("B23").Value = 8.999500000000005

$ws.Range("A28").Value = -22.12729999999999

$ws.Range("A29").Value = -21.28309999999997

$ws.Range("A32").Value = -21.12769999999999

$ws.Range("B34").Value = 9.931600000000008

$ws.Range("E36").Value = 16.39050000000001

$ws.Range("A40").Value = -20.27709999999999

$ws.Range("B43").Value = 5.855

$ws.Range("B45").Value = 4.894000000000002

$ws.Range("E46").Value = 17.03799999999999

$ws.Range("B50").Value = 5.083099999999996
$ws.Range("E50").Value = 16.2823

$ws.Range("B51").Value = 5.691799999999997

$ws.Range("A52").Value = -22.10809999999999

$ws.Range("A57").Value = -22.42230000000001

$ws.Range("A66").Value = -21.4306
$ws.Range("B66").Value = 5.754199999999996

$ws.Range("B67").Value = 5.318499999999998

$ws.Range("B79").Value = 10.02830000000001

$ws.Range("B84").Value = 5.356899999999999

$ws.Range("B92").Value = 4.598599999999998

$ws.Range("E95").Value = 17.72680000000001

$ws.Range("B97").Value = 5.624799999999997
$ws.Range("E97").Value = 16.54449999999999

$ws.Range("A100").Value = -22.0247
